$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultipleChoice")

# Row 5 - "What animal barks?"
$ws.Range("A5").Value = "What animal barks?"
$ws.Range("B5").Value = "Cat"
$ws.Range("C5").Value = "Snake"
$ws.Range("D5").Value = "Pig"
$ws.Range("E5").Value = "Dog"
$ws.Range("F5").Value = "All of the Above"
$ws.Range("G5").Value = "d"

# Row 6 - "Where is Chao Phara River"
$ws.Range("A6").Value = "Where is Chao Phara River"
$ws.Range("B6").Value = "Bkk"
$ws.Range("C6").Value = "Russia"
$ws.Range("D6").Value = "England"
$ws.Range("E6").Value = "Germany"
$ws.Range("F6").Value = "Your Home"
$ws.Range("G6").Value = "a"

# Images for rows 5 and 6
$ws.Range("H5").Value = "gg.png"
$ws.Range("H6").Value = "kk.png"

# Row 7 - "What is Eiffle Tower?" (B7 filled in last, below)
$ws.Range("A7").Value = "What is Eiffle Tower?"
$ws.Range("C7").Value = "Your Grandpa"
$ws.Range("D7").Value = "Egypt"
$ws.Range("E7").Value = "Nitrosomonous"
$ws.Range("F7").Value = "Nitrobacter"
$ws.Range("G7").Value = "a"
$ws.Range("H7").Value = "jj.png"

# Row 8 category first (introduces the new "food" category string)
$ws.Range("I8").Value = "food"

# Row 8 - "What is the worst fruit ever?"
$ws.Range("A8").Value = "What is the worst fruit ever?"
$ws.Range("B8").Value = "Grapes"
$ws.Range("C8").Value = "Apple"
$ws.Range("D8").Value = "Durian"
$ws.Range("E8").Value = "Mango"
$ws.Range("F8").Value = "Starfish"
$ws.Range("G8").Value = "c"
$ws.Range("H8").Value = "hehe.png"

# Re-categorise the existing "Which of the following is a fruit?" question (row 2) as food
$ws.Range("I2").Value = "food"

# Categories for the new rows
$ws.Range("I5").Value = "animal"
$ws.Range("I6").Value = "places"
$ws.Range("I7").Value = "places"

# Finally come back and fill in B7 (the last thing the author touched)
$ws.Range("B7").Value = "BigTower"

# The author ends up with MultipleChoice selected (instead of Matching) with B7 highlighted
$ws.Activate()
$ws.Range("B7").Select()
